# Weekly update: insert a new daily price record for "Repollo" (Macroferia
# Regional de Talca) as row 150, shifting the existing rows 150-190 down to
# 151-191.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 150; everything below (old rows
# 150-190) shifts down to 151-191, preserving all of their data/styles.
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with the new record.
$ws.Range("A150").Value = 5
$ws.Range("B150").Value = "Macroferia Regional de Talca"
$ws.Range("C150").Value = "Maule"
$ws.Range("D150").Value = 44508
$ws.Range("E150").Value = 7
$ws.Range("F150").Value = 100112006
$ws.Range("G150").Value = "Repollo"
$ws.Range("H150").Value = "Crespo record"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 6000
$ws.Range("K150").Value = 700
$ws.Range("L150").Value = 700
$ws.Range("M150").Value = 700
$ws.Range("N150").Value = "$/unidad"
$ws.Range("O150").Value = "Provincia del Elquí"
$ws.Range("P150").Value = 700
$ws.Range("Q150").Value = 1
$ws.Range("R150").Value = "Hortaliza"
